$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.167948961257935
$ws.Range("B1").Value = 2.360015392303467
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.329549312591553
$ws.Range("E1").Value = 1.233513832092285
